# edit.ps1 - applies the "final update - passed exam!" commit to LogSearcher.pptx
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the auto-date placeholder text on the slide master and every
#    slide layout from "15-01-2023" to "16-01-2023".
# ---------------------------------------------------------------------------
function Update-DateShape($shp) {
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "15-01-2023") {
            $tr.Text = "16-01-2023"
        }
    }
}

$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    Update-DateShape $master.Shapes.Item($j)
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $lay = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $lay.Shapes.Count; $j++) {
        Update-DateShape $lay.Shapes.Item($j)
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 1 (title slide): widen/recenter the title box and change the
#    title text from "LogSearcher" to "Udvikling af LogSearcher".
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleShp = $slide1.Shapes.Item(1)

$titleShp.Left = 994299.0 / 914400.0 * 72.0
$titleShp.Width = 8868792.0 / 914400.0 * 72.0

$titleShp.TextFrame.TextRange.Text = "Udvikling af LogSearcher"

# ---------------------------------------------------------------------------
# 3) Slide at position 5 ("LogSearcher komponenter"): remove the
#    "(C# Blazor)" aside from the nested rectangle shape, leaving the line
#    break that precedes it intact.
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$grp4 = $slide5.Shapes.Item(2)
$rect1 = $grp4.GroupItems.Item(1)
$rectTr = $rect1.TextFrame.TextRange
$idx = $rectTr.Text.IndexOf("(C# Blazor)")
if ($idx -ge 0) {
    $delRange = $rectTr.Characters($idx + 1, "(C# Blazor)".Length)
    $delRange.Delete()
}

# ---------------------------------------------------------------------------
# 4) Slide 9 ("Graphtyper - Binary Trees"): widen the explanatory textbox,
#    break it into two paragraphs and bold "kun" and "2".
# ---------------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$textShp = $slide9.Shapes.Item(29)

$textShp.Width = 5075833.0 / 914400.0 * 72.0

$textTr = $textShp.TextFrame.TextRange
$textTr.Text = "Et binary tree karakteriseres ved at `rhver node kun har 2 undernoder (children)"

$kunRange = $textTr.Characters(49, 3)
$kunRange.Font.Bold = 1

$twoRange = $textTr.Characters(57, 1)
$twoRange.Font.Bold = 1

# ---------------------------------------------------------------------------
# 5) Reorder: swap the two slides currently at position 4 and position 5
#    ("kor script 1resetG3.cmd..." and "LogSearcher komponenter").
# ---------------------------------------------------------------------------
$p.Slides.Item(5).MoveTo(4)

Write-Host "edit complete"
